# "added new lists downloads"
# Rename two bull/semen list entries and bump the printed list date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A11 : "Madison" -> "Marquis"
$ws.Range("A11").Value = "Marquis"

# A40 : "Silian-9" -> "Silian-10"
$ws.Range("A40").Value = "Silian-10"

# C1 : list date 19.11.2021 -> 03.12.2021 (serial 44519 -> 44533)
$ws.Range("C1").Value = 44533

# Cosmetic: scroll position / selection as left by the author, and a
# slightly narrower column C.
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("H28").Select()
$ws.Columns.Item(3).ColumnWidth = 37.83
